# On-Orbital-Mechanics.docx edit: "Scale, separate systems and sprite update"
#
# Appends three new paragraphs (the "#5" progress note) to the very end of
# the document body, after the final drawing/paragraph and before the
# section properties.

$d = $word.ActiveDocument

# Start from the end of the document content and collapse to an insertion
# point there, so subsequent InsertParagraphAfter calls append new
# paragraphs rather than editing existing ones.
$end = $d.Content
$end.Collapse(0)

# --- New paragraph 1 -------------------------------------------------
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Text = "#5: I have added basic sprite work to the planets. I have decided that although the sprite work will go unnoticed on the scale of the universe, I want to create a separate scene for each individual planetary system. This scene will have a larger smaller ratio so that the planets and their satellites can be more easily observed."

# --- New paragraph 2 -------------------------------------------------
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Text = "System scales are now separated in their own Galactic scale classes to accommodate for different scales in different systems."

# --- New paragraph 3 -------------------------------------------------
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Text = "I have also corrected the orbital velocity of each planet, using the Semi-Major axis of each around the sun to calculate the appropriate angular velocity. This far more accurately reflects the actual speed they are travelling at any one time, and isn’t obviously wrong by simply plugging in the orbital speed into the angular speed."
